# Auto-generated script applying numeric updates to the Leve profit tables
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 339.33334  # ALC!H2
$ws.Cells.Item(2, 9).Value = 210.53847  # ALC!I2
$ws.Cells.Item(2, 11).Value = 210.53847  # ALC!K2
$ws.Cells.Item(2, 13).Value = -97.53846999999999  # ALC!M2

$ws.Cells.Item(6, 8).Value = 906.4545000000001  # ALC!H6
$ws.Cells.Item(6, 9).Value = 121.375  # ALC!I6
$ws.Cells.Item(6, 10).Value = 3000  # ALC!J6
$ws.Cells.Item(6, 11).Value = 364.125  # ALC!K6
$ws.Cells.Item(6, 12).Value = 9000  # ALC!L6
$ws.Cells.Item(6, 13).Value = -252.125  # ALC!M6
$ws.Cells.Item(6, 14).Value = -9224  # ALC!N6

$ws.Cells.Item(9, 8).Value = 104  # ALC!H9
$ws.Cells.Item(9, 9).Value = 120.6  # ALC!I9
$ws.Cells.Item(9, 11).Value = 120.6  # ALC!K9
$ws.Cells.Item(9, 13).Value = 48.40000000000001  # ALC!M9

$ws.Cells.Item(26, 8).Value = 9999  # ALC!H26
$ws.Cells.Item(26, 10).Value = 9999  # ALC!J26
$ws.Cells.Item(26, 12).Value = 9999  # ALC!L26
$ws.Cells.Item(26, 14).Value = -10687  # ALC!N26

$ws.Cells.Item(76, 8).Value = 0  # ALC!H76
$ws.Cells.Item(76, 9).Value = 0  # ALC!I76
$ws.Cells.Item(76, 11).Value = 0  # ALC!K76
$ws.Cells.Item(76, 13).Value = $null  # ALC!M76 (was -7685)

$ws.Cells.Item(79, 8).Value = 0  # ALC!H79
$ws.Cells.Item(79, 9).Value = 0  # ALC!I79
$ws.Cells.Item(79, 11).Value = 0  # ALC!K79
$ws.Cells.Item(79, 13).Value = $null  # ALC!M79 (was -6908)

$ws.Cells.Item(98, 8).Value = 834.0769  # ALC!H98
$ws.Cells.Item(98, 9).Value = 834.0769  # ALC!I98
$ws.Cells.Item(98, 11).Value = 834.0769  # ALC!K98
$ws.Cells.Item(98, 13).Value = 663.9231  # ALC!M98

$ws.Cells.Item(113, 8).Value = 1342.3572  # ALC!H113
$ws.Cells.Item(113, 9).Value = 1291.8462  # ALC!I113
$ws.Cells.Item(113, 11).Value = 1291.8462  # ALC!K113
$ws.Cells.Item(113, 13).Value = 1962.1538  # ALC!M113

$ws.Cells.Item(122, 8).Value = 834.0769  # ALC!H122
$ws.Cells.Item(122, 9).Value = 834.0769  # ALC!I122
$ws.Cells.Item(122, 11).Value = 2502.2307  # ALC!K122
$ws.Cells.Item(122, 13).Value = -52.23070000000007  # ALC!M122

$ws.Cells.Item(125, 8).Value = 3332  # ALC!H125
$ws.Cells.Item(125, 9).Value = 2498  # ALC!I125
$ws.Cells.Item(125, 11).Value = 22482  # ALC!K125
$ws.Cells.Item(125, 13).Value = -20022  # ALC!M125

$ws.Cells.Item(131, 8).Value = 2204.4  # ALC!H131
$ws.Cells.Item(131, 9).Value = 2362.5  # ALC!I131
$ws.Cells.Item(131, 10).Value = 2099  # ALC!J131
$ws.Cells.Item(131, 11).Value = 7087.5  # ALC!K131
$ws.Cells.Item(131, 12).Value = 6297  # ALC!L131
$ws.Cells.Item(131, 13).Value = -2047.5  # ALC!M131
$ws.Cells.Item(131, 14).Value = -16377  # ALC!N131

$ws.Cells.Item(135, 8).Value = 3057.75  # ALC!H135
$ws.Cells.Item(135, 10).Value = 3166  # ALC!J135
$ws.Cells.Item(135, 12).Value = 28494  # ALC!L135
$ws.Cells.Item(135, 14).Value = -33564  # ALC!N135

$ws.Cells.Item(137, 8).Value = 3150.2593  # ALC!H137
$ws.Cells.Item(137, 9).Value = 1982.6666  # ALC!I137
$ws.Cells.Item(137, 11).Value = 5947.9998  # ALC!K137
$ws.Cells.Item(137, 13).Value = -3397.9998  # ALC!M137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3263.5386  # ARM!H61
$ws.Cells.Item(61, 9).Value = 2373.9524  # ARM!I61
$ws.Cells.Item(61, 11).Value = 2373.9524  # ARM!K61
$ws.Cells.Item(61, 13).Value = -2161.9524  # ARM!M61

$ws.Cells.Item(74, 8).Value = 1490.7142  # ARM!H74
$ws.Cells.Item(74, 9).Value = 1506.3636  # ARM!I74
$ws.Cells.Item(74, 10).Value = 1433.3334  # ARM!J74
$ws.Cells.Item(74, 11).Value = 1506.3636  # ARM!K74
$ws.Cells.Item(74, 12).Value = 1433.3334  # ARM!L74
$ws.Cells.Item(74, 13).Value = -632.3635999999999  # ARM!M74
$ws.Cells.Item(74, 14).Value = -3181.3334  # ARM!N74

$ws.Cells.Item(77, 8).Value = 1490.7142  # ARM!H77
$ws.Cells.Item(77, 9).Value = 1506.3636  # ARM!I77
$ws.Cells.Item(77, 10).Value = 1433.3334  # ARM!J77
$ws.Cells.Item(77, 11).Value = 7531.817999999999  # ARM!K77
$ws.Cells.Item(77, 12).Value = 7166.666999999999  # ARM!L77
$ws.Cells.Item(77, 13).Value = -3163.817999999999  # ARM!M77
$ws.Cells.Item(77, 14).Value = -15902.667  # ARM!N77

$ws.Cells.Item(110, 8).Value = 3725  # ARM!H110
$ws.Cells.Item(110, 9).Value = 1700  # ARM!I110
$ws.Cells.Item(110, 11).Value = 1700  # ARM!K110
$ws.Cells.Item(110, 13).Value = 345  # ARM!M110

$ws.Cells.Item(122, 8).Value = 1684.2858  # ARM!H122
$ws.Cells.Item(122, 9).Value = 1915.1666  # ARM!I122
$ws.Cells.Item(122, 10).Value = 299  # ARM!J122
$ws.Cells.Item(122, 11).Value = 5745.4998  # ARM!K122
$ws.Cells.Item(122, 12).Value = 897  # ARM!L122
$ws.Cells.Item(122, 13).Value = -3295.4998  # ARM!M122
$ws.Cells.Item(122, 14).Value = -5797  # ARM!N122

$ws.Cells.Item(132, 8).Value = 3769.1538  # ARM!H132
$ws.Cells.Item(132, 9).Value = 2100  # ARM!I132
$ws.Cells.Item(132, 11).Value = 6300  # ARM!K132
$ws.Cells.Item(132, 13).Value = -3770  # ARM!M132

$ws.Cells.Item(136, 8).Value = 3263.5386  # ARM!H136
$ws.Cells.Item(136, 9).Value = 2373.9524  # ARM!I136
$ws.Cells.Item(136, 11).Value = 7121.8572  # ARM!K136
$ws.Cells.Item(136, 13).Value = -4571.8572  # ARM!M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 6700.8237  # BSM!H86
$ws.Cells.Item(86, 9).Value = 5496.4  # BSM!I86
$ws.Cells.Item(86, 11).Value = 5496.4  # BSM!K86
$ws.Cells.Item(86, 13).Value = -4373.4  # BSM!M86

$ws.Cells.Item(89, 8).Value = 6700.8237  # BSM!H89
$ws.Cells.Item(89, 9).Value = 5496.4  # BSM!I89
$ws.Cells.Item(89, 11).Value = 27482  # BSM!K89
$ws.Cells.Item(89, 13).Value = -21866  # BSM!M89

$ws.Cells.Item(94, 8).Value = 621.6923  # BSM!H94
$ws.Cells.Item(94, 10).Value = 700.6667  # BSM!J94
$ws.Cells.Item(94, 12).Value = 700.6667  # BSM!L94
$ws.Cells.Item(94, 14).Value = -1602.6667  # BSM!N94

$ws.Cells.Item(134, 8).Value = 4300.4614  # BSM!H134
$ws.Cells.Item(134, 9).Value = 4072.48  # BSM!I134
$ws.Cells.Item(134, 11).Value = 12217.44  # BSM!K134
$ws.Cells.Item(134, 13).Value = -9682.440000000001  # BSM!M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 4808.778  # CRP!H58
$ws.Cells.Item(58, 9).Value = 2558.95  # CRP!I58
$ws.Cells.Item(58, 11).Value = 2558.95  # CRP!K58
$ws.Cells.Item(58, 13).Value = -2355.95  # CRP!M58

$ws.Cells.Item(68, 8).Value = 72260.336  # CRP!H68
$ws.Cells.Item(68, 10).Value = 72260.336  # CRP!J68
$ws.Cells.Item(68, 12).Value = 72260.336  # CRP!L68
$ws.Cells.Item(68, 14).Value = -73758.336  # CRP!N68

$ws.Cells.Item(71, 8).Value = 72260.336  # CRP!H71
$ws.Cells.Item(71, 10).Value = 72260.336  # CRP!J71
$ws.Cells.Item(71, 12).Value = 216781.008  # CRP!L71
$ws.Cells.Item(71, 14).Value = -224269.008  # CRP!N71

$ws.Cells.Item(96, 8).Value = 12516.454  # CRP!H96
$ws.Cells.Item(96, 10).Value = 12516.454  # CRP!J96
$ws.Cells.Item(96, 12).Value = 12516.454  # CRP!L96
$ws.Cells.Item(96, 14).Value = -18008.454  # CRP!N96

$ws.Cells.Item(132, 8).Value = 7247.5  # CRP!H132
$ws.Cells.Item(132, 9).Value = 4500  # CRP!I132
$ws.Cells.Item(132, 11).Value = 13500  # CRP!K132
$ws.Cells.Item(132, 13).Value = -10970  # CRP!M132

$ws.Cells.Item(136, 8).Value = 4808.778  # CRP!H136
$ws.Cells.Item(136, 9).Value = 2558.95  # CRP!I136
$ws.Cells.Item(136, 11).Value = 7676.849999999999  # CRP!K136
$ws.Cells.Item(136, 13).Value = -5126.849999999999  # CRP!M136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(57, 8).Value = 1472.5  # CUL!H57
$ws.Cells.Item(57, 10).Value = 2965  # CUL!J57
$ws.Cells.Item(57, 12).Value = 8895  # CUL!L57
$ws.Cells.Item(57, 14).Value = -10013  # CUL!N57

$ws.Cells.Item(60, 8).Value = 813.6429000000001  # CUL!H60
$ws.Cells.Item(60, 9).Value = 298.2143  # CUL!I60
$ws.Cells.Item(60, 10).Value = 1329.0714  # CUL!J60
$ws.Cells.Item(60, 11).Value = 894.6428999999999  # CUL!K60
$ws.Cells.Item(60, 12).Value = 3987.2142  # CUL!L60
$ws.Cells.Item(60, 13).Value = -643.6428999999999  # CUL!M60
$ws.Cells.Item(60, 14).Value = -4489.2142  # CUL!N60

$ws.Cells.Item(131, 10).Value = 1000  # CUL!J131
$ws.Cells.Item(131, 12).Value = 3000  # CUL!L131
$ws.Cells.Item(131, 14).Value = -13080  # CUL!N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4008.5  # GSM!H70
$ws.Cells.Item(70, 9).Value = 3619.9092  # GSM!I70
$ws.Cells.Item(70, 10).Value = 5433.3335  # GSM!J70
$ws.Cells.Item(70, 11).Value = 3619.9092  # GSM!K70
$ws.Cells.Item(70, 12).Value = 5433.3335  # GSM!L70
$ws.Cells.Item(70, 13).Value = -3349.9092  # GSM!M70
$ws.Cells.Item(70, 14).Value = -5973.3335  # GSM!N70

$ws.Cells.Item(73, 8).Value = 4008.5  # GSM!H73
$ws.Cells.Item(73, 9).Value = 3619.9092  # GSM!I73
$ws.Cells.Item(73, 10).Value = 5433.3335  # GSM!J73
$ws.Cells.Item(73, 11).Value = 3619.9092  # GSM!K73
$ws.Cells.Item(73, 12).Value = 5433.3335  # GSM!L73
$ws.Cells.Item(73, 13).Value = -2683.9092  # GSM!M73
$ws.Cells.Item(73, 14).Value = -7305.3335  # GSM!N73

$ws.Cells.Item(97, 8).Value = 921.2  # GSM!H97
$ws.Cells.Item(97, 9).Value = 921.2  # GSM!I97
$ws.Cells.Item(97, 10).Value = 0  # GSM!J97
$ws.Cells.Item(97, 11).Value = 921.2  # GSM!K97
$ws.Cells.Item(97, 12).Value = 0  # GSM!L97
$ws.Cells.Item(97, 13).Value = -425.2  # GSM!M97
$ws.Cells.Item(97, 14).Value = $null  # GSM!N97 (was -2239.5)

$ws.Cells.Item(102, 8).Value = 1957.5  # GSM!H102
$ws.Cells.Item(102, 9).Value = 1592.44  # GSM!I102
$ws.Cells.Item(102, 11).Value = 1592.44  # GSM!K102
$ws.Cells.Item(102, 13).Value = 29.55999999999995  # GSM!M102

$ws.Cells.Item(132, 8).Value = 38802.266  # GSM!H132
$ws.Cells.Item(132, 9).Value = 45209  # GSM!I132
$ws.Cells.Item(132, 11).Value = 135627  # GSM!K132
$ws.Cells.Item(132, 13).Value = -133097  # GSM!M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 2722.75  # LTW!H122
$ws.Cells.Item(122, 9).Value = 2722.75  # LTW!I122
$ws.Cells.Item(122, 11).Value = 8168.25  # LTW!K122
$ws.Cells.Item(122, 13).Value = -5718.25  # LTW!M122

$ws.Cells.Item(132, 8).Value = 3361.5  # LTW!H132
$ws.Cells.Item(132, 9).Value = 3361.5  # LTW!I132
$ws.Cells.Item(132, 11).Value = 10084.5  # LTW!K132
$ws.Cells.Item(132, 13).Value = -7554.5  # LTW!M132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 10000  # WVR!H62
$ws.Cells.Item(62, 10).Value = 12000  # WVR!J62
$ws.Cells.Item(62, 12).Value = 12000  # WVR!L62
$ws.Cells.Item(62, 14).Value = -13248  # WVR!N62

$ws.Cells.Item(65, 8).Value = 10000  # WVR!H65
$ws.Cells.Item(65, 10).Value = 12000  # WVR!J65
$ws.Cells.Item(65, 12).Value = 60000  # WVR!L65
$ws.Cells.Item(65, 14).Value = -66240  # WVR!N65

$ws.Cells.Item(113, 8).Value = 445.52942  # WVR!H113
$ws.Cells.Item(113, 9).Value = 306.66666  # WVR!I113
$ws.Cells.Item(113, 11).Value = 919.9999799999999  # WVR!K113
$ws.Cells.Item(113, 13).Value = 1250.00002  # WVR!M113
